# Update the "Stock Group" worksheet manifest with new DOM stock groups
# (Atmosphere, DOM: Deadwood, DOM: Down Deadwood, DOM: Litter, DOM: Soil,
#  DOM: Standing Deadwood) inserted in their correct alphabetical spots.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Stock Group")

# Final target layout (row -> Name / Description). Rows are written from the
# bottom upward so that source data is never overwritten before it has been
# relocated to its new row, and cell-by-cell assignment is used instead of
# Rows.Insert() so Excel does not fabricate extra formatted-but-empty cells
# in columns B:D (Insert() copies the whole row's formatting across all used
# columns).
$rows = @(
    @("Atmosphere", $null),
    @("Biomass: Aboveground", "Carbon in all aboveground biomass pools"),
    @("Biomass: Belowground", "Carbon in all belowground biomnass pools (coarse plus fine roots)"),
    @("Biomass: Total", "Carbon in aboverground and belowground biomass pools"),
    @("DOM: Aboveground", "Carbon in DOM pools above the mineral soil"),
    @("DOM: Belowground", "Carbon in DOM pools in the mineral soil"),
    @("DOM: Deadwood", "Carbon in all deadwood pools"),
    @("DOM: Down Deadwood", "Carbon stores in down deadwood"),
    @("DOM: Litter", $null),
    @("DOM: Soil", $null),
    @("DOM: Standing Deadwood", "Carbon stored in standing deadwood"),
    @("DOM: Total", "Carbon in all DOM pools"),
    @("Total Ecosystem", "Carbon in Biomass and DOM pools")
)

$firstDataRow = 2

for ($i = $rows.Count - 1; $i -ge 0; $i--) {
    $r = $firstDataRow + $i
    $name = $rows[$i][0]
    $desc = $rows[$i][1]

    $ws.Cells.Item($r, 1).Value = $name
    if ($desc) {
        $ws.Cells.Item($r, 2).Value = $desc
    } else {
        $ws.Cells.Item($r, 2).ClearContents() | Out-Null
        $ws.Cells.Item($r, 2).Style = "Normal"
    }
}

# Column D was hidden (used only for the IsAuto validation helper column) -
# now shown again at its existing best-fit width.
$ws.Columns.Item(4).Hidden = $false

# Match the saved selection state recorded for this sheet.
$ws.Range("B23").Select() | Out-Null
